$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 40271
$ws.Range("J63").Value = 40271
$ws.Range("L63").Value = 40271
$ws.Range("N63").Value = -41519
$ws.Range("H66").Value = 40271
$ws.Range("J66").Value = 40271
$ws.Range("L66").Value = 120813
$ws.Range("N66").Value = -127053
$ws.Range("H132").Value = 1883.5927
$ws.Range("I132").Value = 1883.5927
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5650.7781
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3120.7781
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 35705.824
$ws.Range("J134").Value = 35705.824
$ws.Range("L134").Value = 35705.824
$ws.Range("N134").Value = -45845.824
$ws.Range("H137").Value = 1439.9623
$ws.Range("I137").Value = 946.76666
$ws.Range("J137").Value = 2083.261
$ws.Range("K137").Value = 2840.29998
$ws.Range("L137").Value = 6249.782999999999
$ws.Range("M137").Value = -290.2999799999998
$ws.Range("N137").Value = -11349.783
$ws.Range("H138").Value = 3995.9556
$ws.Range("I138").Value = 1808
$ws.Range("J138").Value = 5324.357
$ws.Range("K138").Value = 5424
$ws.Range("L138").Value = 15973.071
$ws.Range("M138").Value = -284
$ws.Range("N138").Value = -26253.071

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5849.1177
$ws.Range("I32").Value = 5375.2144
$ws.Range("J32").Value = 8060.6665
$ws.Range("K32").Value = 5375.2144
$ws.Range("L32").Value = 8060.6665
$ws.Range("M32").Value = -5088.2144
$ws.Range("N32").Value = -8634.666499999999
$ws.Range("H132").Value = 2201.5088
$ws.Range("I132").Value = 1274.4595
$ws.Range("J132").Value = 3916.55
$ws.Range("K132").Value = 3823.3785
$ws.Range("L132").Value = 11749.65
$ws.Range("M132").Value = -1293.3785
$ws.Range("N132").Value = -16809.65
$ws.Range("H133").Value = 41252.75
$ws.Range("J133").Value = 41252.75
$ws.Range("L133").Value = 41252.75
$ws.Range("N133").Value = -46312.75
$ws.Range("H134").Value = 18950
$ws.Range("J134").Value = 18950
$ws.Range("L134").Value = 18950
$ws.Range("N134").Value = -29090
$ws.Range("H139").Value = 45503.75
$ws.Range("J139").Value = 45503.75
$ws.Range("L139").Value = 45503.75
$ws.Range("N139").Value = -55783.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4244.2666
$ws.Range("I134").Value = 5743.2085
$ws.Range("J134").Value = 2531.1904
$ws.Range("K134").Value = 17229.6255
$ws.Range("L134").Value = 7593.5712
$ws.Range("M134").Value = -14694.6255
$ws.Range("N134").Value = -12663.5712
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H139").Value = 26709
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2476.3247
$ws.Range("I31").Value = 1552.6
$ws.Range("J31").Value = 2800.4385
$ws.Range("K31").Value = 1552.6
$ws.Range("L31").Value = 2800.4385
$ws.Range("M31").Value = -1257.6
$ws.Range("N31").Value = -3390.4385
$ws.Range("H34").Value = 2476.3247
$ws.Range("I34").Value = 1552.6
$ws.Range("J34").Value = 2800.4385
$ws.Range("K34").Value = 1552.6
$ws.Range("L34").Value = 2800.4385
$ws.Range("M34").Value = -1350.6
$ws.Range("N34").Value = -3204.4385
$ws.Range("H135").Value = 34120
$ws.Range("J135").Value = 34120
$ws.Range("L135").Value = 34120
$ws.Range("N135").Value = -44260

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2616.1604
$ws.Range("I68").Value = 3793.7058
$ws.Range("J68").Value = 1764.3191
$ws.Range("K68").Value = 11381.1174
$ws.Range("L68").Value = 5292.9573
$ws.Range("M68").Value = -10570.1174
$ws.Range("N68").Value = -6914.9573
$ws.Range("H69").Value = 500
$ws.Range("I69").Value = 500
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 1500
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -689
$ws.Range("N69").ClearContents()
$ws.Range("H71").Value = 2616.1604
$ws.Range("I71").Value = 3793.7058
$ws.Range("J71").Value = 1764.3191
$ws.Range("K71").Value = 34143.3522
$ws.Range("L71").Value = 15878.8719
$ws.Range("M71").Value = -30087.3522
$ws.Range("N71").Value = -23990.8719
$ws.Range("H72").Value = 500
$ws.Range("I72").Value = 500
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 4500
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -444
$ws.Range("N72").ClearContents()
$ws.Range("H80").Value = 4515.8887
$ws.Range("I80").Value = 1850
$ws.Range("J80").Value = 5277.5713
$ws.Range("K80").Value = 5550
$ws.Range("L80").Value = 15832.7139
$ws.Range("M80").Value = -4614
$ws.Range("N80").Value = -17704.7139
$ws.Range("H83").Value = 4515.8887
$ws.Range("I83").Value = 1850
$ws.Range("J83").Value = 5277.5713
$ws.Range("K83").Value = 16650
$ws.Range("L83").Value = 47498.14169999999
$ws.Range("M83").Value = -11970
$ws.Range("N83").Value = -56858.14169999999
$ws.Range("H86").Value = 1252.75
$ws.Range("I86").Value = 1302.7333
$ws.Range("J86").Value = 503
$ws.Range("K86").Value = 3908.199900000001
$ws.Range("L86").Value = 1509
$ws.Range("M86").Value = -2722.199900000001
$ws.Range("N86").Value = -3881
$ws.Range("H89").Value = 1252.75
$ws.Range("I89").Value = 1302.7333
$ws.Range("J89").Value = 503
$ws.Range("K89").Value = 11724.5997
$ws.Range("L89").Value = 4527
$ws.Range("M89").Value = -5796.599700000001
$ws.Range("N89").Value = -16383
$ws.Range("H107").Value = 1024.4462
$ws.Range("J107").Value = 1234.96
$ws.Range("L107").Value = 3704.88
$ws.Range("N107").Value = -7544.88
$ws.Range("H113").Value = 1304858.1
$ws.Range("I113").Value = 1923549.1
$ws.Range("J113").Value = 500559.9
$ws.Range("K113").Value = 5770647.300000001
$ws.Range("L113").Value = 1501679.7
$ws.Range("M113").Value = -5768477.300000001
$ws.Range("N113").Value = -1506019.7
$ws.Range("H131").Value = 13415800
$ws.Range("I131").Value = 5556038.5
$ws.Range("J131").Value = 15626358
$ws.Range("K131").Value = 16668115.5
$ws.Range("L131").Value = 46879074
$ws.Range("M131").Value = -16663075.5
$ws.Range("N131").Value = -46889154

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2696.1755
$ws.Range("I132").Value = 1815.3334
$ws.Range("J132").Value = 3210
$ws.Range("K132").Value = 5446.0002
$ws.Range("L132").Value = 9630
$ws.Range("M132").Value = -2916.0002
$ws.Range("N132").Value = -14690
$ws.Range("H134").Value = 18392.5
$ws.Range("J134").Value = 18392.5
$ws.Range("L134").Value = 55177.5
$ws.Range("N134").Value = -60247.5
$ws.Range("H135").Value = 35212.668
$ws.Range("J135").Value = 35212.668
$ws.Range("L135").Value = 35212.668
$ws.Range("N135").Value = -45352.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 39626
$ws.Range("J108").Value = 39626
$ws.Range("L108").Value = 39626
$ws.Range("N108").Value = -47306
$ws.Range("H122").Value = 3137987.2
$ws.Range("I122").Value = 4206722.5
$ws.Range("J122").Value = 1119265
$ws.Range("K122").Value = 12620167.5
$ws.Range("L122").Value = 3357795
$ws.Range("M122").Value = -12617717.5
$ws.Range("N122").Value = -3362695
$ws.Range("H132").Value = 13548034
$ws.Range("I132").Value = 28899362
$ws.Range("K132").Value = 86698086
$ws.Range("M132").Value = -86695556
$ws.Range("H134").Value = 42429
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 42429
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 42429
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -52569
$ws.Range("H135").Value = 32500
$ws.Range("I135").Value = 40000
$ws.Range("J135").Value = 25000
$ws.Range("K135").Value = 40000
$ws.Range("L135").Value = 25000
$ws.Range("M135").Value = -34930
$ws.Range("N135").Value = -35140
$ws.Range("H136").Value = 5289.3687
$ws.Range("I136").Value = 4165.814
$ws.Range("J136").Value = 8740.286
$ws.Range("K136").Value = 12497.442
$ws.Range("L136").Value = 26220.858
$ws.Range("M136").Value = -9947.442000000001
$ws.Range("N136").Value = -31320.858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 35423.5
$ws.Range("J46").Value = 35423.5
$ws.Range("L46").Value = 35423.5
$ws.Range("N46").Value = -35885.5
$ws.Range("H126").Value = 907
$ws.Range("I126").Value = 656.86664
$ws.Range("J126").Value = 1376
$ws.Range("K126").Value = 1970.59992
$ws.Range("L126").Value = 4128
$ws.Range("M126").Value = 499.4000800000001
$ws.Range("N126").Value = -9068
$ws.Range("H132").Value = 1444.2029
$ws.Range("I132").Value = 981.04254
$ws.Range("J132").Value = 2433.682
$ws.Range("K132").Value = 2943.12762
$ws.Range("L132").Value = 7301.045999999999
$ws.Range("M132").Value = -413.1276200000002
$ws.Range("N132").Value = -12361.046
$ws.Range("H134").Value = 35423.5
$ws.Range("J134").Value = 35423.5
$ws.Range("L134").Value = 106270.5
$ws.Range("N134").Value = -111340.5
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
